$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{ Row = 2; Old = "Mostly"; New = "Design: Mostly" },
    @{ Row = 3; Old = "Nothing special to mention"; New = "Design: Nothing special to mention" },
    @{ Row = 4; Old = "91xx Went well, some minor budget challenges"; New = "Design: 91xx Went well, some minor budget challenges" },
    @{ Row = 5; Old = "Mostly ok."; New = "Design: Mostly ok." },
    @{ Row = 6; Old = "Internal communication ok. External communication with suppliers mostly ok."; New = "Design: Internal communication ok. External communication with suppliers mostly ok." },
    @{ Row = 7; Old = "Some things went to correct direction but regarding TK the opposite way."; New = "Design: Some things went to correct direction but regarding TK the opposite way." }
)

foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, 2)
    $rng = $cell.Range
    $rng.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 1)
}
